$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.25971245765686
$ws.Range("B1").Value = 2.348007440567017
$ws.Range("C1").Value = 3.522169589996338
$ws.Range("D1").Value = 2.662193059921265
$ws.Range("E1").Value = 1.356282591819763
